$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-21 20:43:56"
$wsZh.Range("H3").Value = "2016-03-21 20:44:20"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-21 20:44:00"
$wsDe.Range("H3").Value = "2016-03-21 20:44:25"
